# fix templates for dictionary uploads (#4940)
# The header of the "Variable values" sheet renamed column C from "name" to
# "variable" (matches the "name"/"variable" columns used on the Variables
# sheet), the column was widened to fit the new header text, and the sheet's
# remembered selection moved to C2.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Variable values")

# C1: "name" -> "variable"
$ws.Range("C1").Value = "variable"

# Column C needs to be a bit wider to fit "variable".
$ws.Columns.Item(3).ColumnWidth = 6

# Remember the new selection (C2) on this sheet, then restore the
# originally active sheet ("Datasets") as the active tab/selection.
$originalActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("C2").Select()
$originalActive.Activate()
